# update test script auto order
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new test case as row 7 (added first so its new "06" shared
# string is registered before the new phone number below, matching the
# order the strings were introduced to the sheet).
$ws.Range("A7").Value = "'06"
$ws.Range("B7").Value = "Ngoc Vu"
$ws.Range("C7").Value = "207 Giai Phong"
$ws.Range("D7").Value = "Ha Noi"
$ws.Range("E7").Value = "'100000"

# Replace the old phone number "0879231543" with the new one "0823456789"
# everywhere it is used (rows 2 and 3 keep a plain text value, row 4 keeps
# the literal leading apostrophe that was already present in the source data).
$ws.Range("F2").Value = "'0823456789"
$ws.Range("F3").Value = "'0823456789"
$ws.Range("F4").Value = "''0823456789"
$ws.Range("F7").Value = "'0823456789"

$ws.Range("G7").Value = "Hmm. We couldn" + [char]8217 + "t find your device" + [char]8217 + "s location. Please enter the address manually."

# Row 7 should wrap like the rest of the data rows.
$ws.Range("A7:G7").WrapText = $true

# B7:D7 pick up an explicit black font color (as opposed to the implicit
# theme color used elsewhere), matching how this row was authored.
$ws.Range("B7:D7").Font.Color = 0

# Row 7 contains a long wrapped message, so it needs extra height.
$ws.Rows.Item(7).RowHeight = 43.2

# Restore the selection to the cell that was active when the sheet was saved.
$null = $ws.Range("E11").Select()
